# Rename the "Data" sheet to "Data table" (per diff: xl/workbook.xml <sheets>)
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")
$dataSheet.Name = "Data table"

# Move the active/selected tab from "Codelists" to the renamed "Data table" sheet
# (per diff: workbookView activeTab 2 -> 3, tabSelected moves from sheet3 to sheet4)
$dataSheet.Activate()
